$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 4666.3335
$ws.Range("J97").Value = 4666.3335
$ws.Range("L97").Value = 13999.0005
$ws.Range("N97").Value = -14991.0005

$ws.Range("H135").Value = 4867.5835
$ws.Range("I135").Value = 801.375
$ws.Range("K135").Value = 7212.375
$ws.Range("M135").Value = -4677.375

$ws.Range("H137").Value = 25003250
$ws.Range("I137").Value = 50001028
$ws.Range("J137").Value = 5471.25
$ws.Range("K137").Value = 150003084
$ws.Range("L137").Value = 16413.75
$ws.Range("M137").Value = -150000534
$ws.Range("N137").Value = -21513.75

$ws.Range("H138").Value = 3219.0244
$ws.Range("I138").Value = 2140.889
$ws.Range("J138").Value = 3522.25
$ws.Range("K138").Value = 6422.667
$ws.Range("L138").Value = 10566.75
$ws.Range("M138").Value = -1282.667
$ws.Range("N138").Value = -20846.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 184750
$ws.Range("I34").Value = 184750
$ws.Range("K34").Value = 184750
$ws.Range("M34").Value = -184479

$ws.Range("H61").Value = 3575806
$ws.Range("I61").Value = 4410.5386
$ws.Range("J61").Value = 50003948
$ws.Range("K61").Value = 4410.5386
$ws.Range("L61").Value = 50003948
$ws.Range("M61").Value = -4198.5386
$ws.Range("N61").Value = -50004372

$ws.Range("H136").Value = 3575806
$ws.Range("I136").Value = 4410.5386
$ws.Range("J136").Value = 50003948
$ws.Range("K136").Value = 13231.6158
$ws.Range("L136").Value = 150011844
$ws.Range("M136").Value = -10681.6158
$ws.Range("N136").Value = -150016944

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2843.0256
$ws.Range("I105").Value = 2426.862
$ws.Range("K105").Value = 2426.862
$ws.Range("M105").Value = -679.8620000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1986776.1
$ws.Range("I31").Value = 2418195.8
$ws.Range("K31").Value = 2418195.8
$ws.Range("M31").Value = -2417900.8

$ws.Range("H34").Value = 1986776.1
$ws.Range("I34").Value = 2418195.8
$ws.Range("K34").Value = 2418195.8
$ws.Range("M34").Value = -2417993.8

$ws.Range("H99").Value = 27939.666
$ws.Range("I99").Value = 28166.555
$ws.Range("J99").Value = 27259
$ws.Range("K99").Value = 28166.555
$ws.Range("L99").Value = 27259
$ws.Range("M99").Value = -26668.555
$ws.Range("N99").Value = -30255

$ws.Range("H126").Value = 27939.666
$ws.Range("I126").Value = 28166.555
$ws.Range("J126").Value = 27259
$ws.Range("K126").Value = 84499.66500000001
$ws.Range("L126").Value = 81777
$ws.Range("M126").Value = -82029.66500000001
$ws.Range("N126").Value = -86717

$ws.Range("H134").Value = 3013.6155
$ws.Range("I134").Value = 2697.9092
$ws.Range("K134").Value = 8093.7276
$ws.Range("M134").Value = -5558.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 250000530
$ws.Range("I11").Value = 1044.5
$ws.Range("J11").Value = 500000030
$ws.Range("K11").Value = 3133.5
$ws.Range("L11").Value = 1500000090
$ws.Range("M11").Value = -2993.5
$ws.Range("N11").Value = -1500000370

$ws.Range("H29").Value = 1026.4
$ws.Range("I29").Value = 1144
$ws.Range("J29").Value = 850
$ws.Range("K29").Value = 3432
$ws.Range("L29").Value = 2550
$ws.Range("M29").Value = -3155
$ws.Range("N29").Value = -3104

$ws.Range("H34").Value = 353
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168

$ws.Range("H52").Value = 3947.5
$ws.Range("J52").Value = 3947.5
$ws.Range("L52").Value = 11842.5
$ws.Range("N52").Value = -12374.5

$ws.Range("H86").Value = 241.57143
$ws.Range("J86").Value = 298.2
$ws.Range("L86").Value = 894.5999999999999
$ws.Range("N86").Value = -3266.6

$ws.Range("H89").Value = 241.57143
$ws.Range("J89").Value = 298.2
$ws.Range("L89").Value = 2683.8
$ws.Range("N89").Value = -14539.8

$ws.Range("H107").Value = 3510.4375
$ws.Range("I107").Value = 579.6
$ws.Range("J107").Value = 4842.636
$ws.Range("K107").Value = 1738.8
$ws.Range("L107").Value = 14527.908
$ws.Range("M107").Value = 181.1999999999998
$ws.Range("N107").Value = -18367.908

$ws.Range("H109").Value = 13000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 13000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 39000
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -41080

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H132").Value = 1988.4445
$ws.Range("I132").Value = 1250
$ws.Range("K132").Value = 11250
$ws.Range("M132").Value = -8720

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2167.375
$ws.Range("I102").Value = 2111.8667
$ws.Range("K102").Value = 2111.8667
$ws.Range("M102").Value = -489.8667

$ws.Range("H132").Value = 5092.4243
$ws.Range("I132").Value = 4524.0586
$ws.Range("K132").Value = 13572.1758
$ws.Range("M132").Value = -11042.1758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2067.7646
$ws.Range("I22").Value = 559.2
$ws.Range("K22").Value = 559.2
$ws.Range("M22").Value = -264.2

$ws.Range("H27").Value = 2067.7646
$ws.Range("I27").Value = 559.2
$ws.Range("K27").Value = 559.2
$ws.Range("M27").Value = -452.2

$ws.Range("H74").Value = 49758
$ws.Range("I74").Value = 49758
$ws.Range("K74").Value = 49758
$ws.Range("M74").Value = -48760

$ws.Range("H77").Value = 49758
$ws.Range("I77").Value = 49758
$ws.Range("K77").Value = 149274
$ws.Range("M77").Value = -144282

$ws.Range("H132").Value = 5083393
$ws.Range("I132").Value = 8348813
$ws.Range("K132").Value = 25046439
$ws.Range("M132").Value = -25043909

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 49000
$ws.Range("J64").Value = 49000
$ws.Range("L64").Value = 49000
$ws.Range("N64").Value = -49496

$ws.Range("H67").Value = 49000
$ws.Range("J67").Value = 49000
$ws.Range("L67").Value = 49000
$ws.Range("N67").Value = -50716

$ws.Range("H113").Value = 922.2
$ws.Range("I113").Value = 980.3684
$ws.Range("J113").Value = 821.7273
$ws.Range("K113").Value = 2941.1052
$ws.Range("L113").Value = 2465.1819
$ws.Range("M113").Value = -771.1052
$ws.Range("N113").Value = -6805.1819

$ws.Range("H132").Value = 3790421
$ws.Range("I132").Value = 4506915.5
$ws.Range("J132").Value = 3236.2856
$ws.Range("K132").Value = 13520746.5
$ws.Range("M132").Value = -13518216.5
$ws.Range("N132").Value = -14768.8568
